# "cleaned the code a bit"
# Applies the recorded edit: restructure the "ERD" sheet into a wide
# header table (rows 1-5) followed by the original single-column field
# list (now rows 8-35), tweak a shared-formula range on "2test", and
# move the active selection around a couple of sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "ERD" sheet (2nd tab): insert 7 rows at the top so the existing
#    A1:F28 field-list content shifts down to A8:F35, then build the
#    new wide header block in rows 1-5, clear the now-orphaned F column
#    duplicate entries, and fix the one value that differs (row 13:
#    "Company Name" -> "Customer Name").
# ---------------------------------------------------------------------
$erd = $wb.Worksheets.Item(2)

$erd.Rows("1:7").Insert()

# Style reference cells (post-shift): these already carry the exact
# cellXf indices the new header rows need to reuse.
$styleBoldRed  = $erd.Range("A8")    # s=3
$styleNormal   = $erd.Range("A9")    # s=2
$styleContact  = $erd.Range("A20")   # s=5
$styleFooter   = $erd.Range("A34")   # s=4

# Row 1
$erd.Range("A1").Value = "Unique Lead Assignment Number "
$erd.Range("B1").Value = "Customer Name"
$erd.Range("C1").Value = "Address Line 1"
$erd.Range("D1").Value = "Address Line 2"
$erd.Range("E1").Value = "City"
$erd.Range("F1").Value = "State"
$erd.Range("G1").Value = "Post Code"
$erd.Range("H1").Value = "Industry "
$erd.Range("I1").Value = "Physical Channel"
$erd.Range("J1").Value = "Main Phone #"
$erd.Range("K1").Value = "Website"
$erd.Range("L1").Value = "SSM No"

# Row 2
$erd.Range("A2").Value = "Unique Lead Assignment Number "
$erd.Range("B2").Value = "Customer Name"
$erd.Range("C2").Value = "Competitors"
$erd.Range("D2").Value = "Total Potential Revenue/Month"

# Row 3
$erd.Range("A3").Value = "Unique Lead Assignment Number "
$erd.Range("B3").Value = "Customer Name"
$erd.Range("C3").Value = "Contact Person Name "
$erd.Range("D3").Value = "Contact Person Email"
$erd.Range("E3").Value = "Contact person Designation"
$erd.Range("F3").Value = "Contact Person Phone"

# Row 4
$erd.Range("A4").Value = "Unique Lead Assignment Number "
$erd.Range("B4").Value = "Customer Name"
$erd.Range("C4").Value = "Lead Originator"
$erd.Range("D4").Value = "Lead Originator Email"
$erd.Range("E4").Value = "Created Date"
$erd.Range("F4").Value = "Created By"
$erd.Range("G4").Value = "Suspect Accepted By"
$erd.Range("H4").Value = "Suspect Accepted At"
$erd.Range("I4").Value = "Prospect Accepted By"
$erd.Range("J4").Value = "Prospect Accepted At"

# Row 5
$erd.Range("A5").Value = "Unique Lead Assignment Number "
$erd.Range("B5").Value = "Customer Name"
$erd.Range("C5").Value = "Source Type"
$erd.Range("D5").Value = "Lead Priority Level"

# Apply the matching cell formatting to every new header cell (value
# already set above; PasteSpecial(formats) only touches the style).
$styleBoldRed.Copy()
$erd.Range("A1,A2,A3,A4,A5").PasteSpecial(-4122)

$styleNormal.Copy()
$erd.Range("B1:L2,B3:B5,C4:J4").PasteSpecial(-4122)

$styleContact.Copy()
$erd.Range("C3:F3").PasteSpecial(-4122)

$styleFooter.Copy()
$erd.Range("C5:D5").PasteSpecial(-4122)

$erd.Application.CutCopyMode = $false

# K4 stays a blank, footer-styled cell (matches <c r="K4" s="4"/>).
$styleFooter.Copy()
$erd.Range("K4").PasteSpecial(-4122)
$erd.Range("K4").ClearContents()

# The old F-column duplicate list (now at F19:F23) is dropped entirely.
$erd.Range("F19:F23").Clear()

# The single value swap in the single-column list: old "Company Name"
# (now row 13) becomes "Customer Name".
$erd.Range("A13").Value = "Customer Name"

# Move the ERD selection/active cell and make this the active tab.
$erd.Activate()
$erd.Range("D21").Select()

# ---------------------------------------------------------------------
# 2) "Sample Records " sheet (3rd tab): move the saved selection.
# ---------------------------------------------------------------------
$sampleRecords = $wb.Worksheets.Item(3)
$sampleRecords.Activate()
$sampleRecords.Range("F1").Select()

# ---------------------------------------------------------------------
# 3) "2test" sheet (5th tab): shrink the shared-formula declared range
#    from F3:F13 to F3:F12 by re-entering it across the real extent.
# ---------------------------------------------------------------------
$twotest = $wb.Worksheets.Item(5)
$twotest.Range("F3:F5").FormulaR1C1 = '=_xlfn.CONCAT(RC[-4],"@",RC[-1])'
$twotest.Range("F10:F12").FormulaR1C1 = '=_xlfn.CONCAT(RC[-4],"@",RC[-1])'

# ---------------------------------------------------------------------
# 4) Leave "ERD" as the active/selected tab last (matches activeTab=1,
#    and drops tabSelected from "fuzzy" which had it before).
# ---------------------------------------------------------------------
$erd.Activate()
$erd.Range("D21").Select()
